$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add description text for the "Writing/Narrative Designer" position.
#    It is the 2nd table in the document, row 2 ("Description:"), 2nd cell,
#    which is currently an empty paragraph.
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(2)
$cell = $tbl.Rows.Item(2).Cells.Item(2)
$rng = $cell.Range
# Collapse to just before the cell's trailing paragraph/end-of-cell marker
$rng.End = $rng.End - 1
$rng.Text = "Focuses on creating a story/plot for every game they develop, focusing on the player`u{2019}s actions and choices so new dialogue is added for testers to discover."
$rng.Font.Bold = 1

# ---------------------------------------------------------------------------
# 2) Merge the header runs around "SuniTAFE" into a single run / text node,
#    removing the spell-check proofing marks split around it.
# ---------------------------------------------------------------------------
$sections = $d.Sections
for ($i = 1; $i -le $sections.Count; $i++) {
  $hf = $sections.Item($i).Headers.Item(1)
  if ($hf.Exists) {
    $hf.Range.Find.Execute(" SuniTAFE ", $false, $false, $false, $false, $false, $true, 1, $false, " SuniTAFE ", 2)
  }
}
